$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# The team-name labels in column B for rows 2/3 and rows 7/8 are swapped
# (shared-string pool re-ordered upstream), and the Weekly Pending Total
# (C) / Repayment (D) figures for every team row are updated. Column E
# is a formula (=D/C) and recalculates automatically.

$ws.Range("B2").Value = "Cpu_s2l"
$ws.Range("B3").Value = "Xinghao_s2l"
$ws.Range("B7").Value = "Zakka_S2l"
$ws.Range("B8").Value = "Mkm_s2l"

$ws.Range("C2").Value = 1705546967
$ws.Range("D2").Value = 179546829

$ws.Range("C3").Value = 2540164470
$ws.Range("D3").Value = 266165525

$ws.Range("C4").Value = 6145029203
$ws.Range("D4").Value = 643182863

$ws.Range("C5").Value = 1689163205
$ws.Range("D5").Value = 164348626

$ws.Range("C6").Value = 6227337575
$ws.Range("D6").Value = 572049849

$ws.Range("C7").Value = 6093487702
$ws.Range("D7").Value = 534233622

$ws.Range("C8").Value = 3295574061
$ws.Range("D8").Value = 287693623

$wb.Application.CalculateFull()
